$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "Girmaa Masqalaa" row (row 5) entirely.
$ws.Rows.Item(5).Delete()

# After the first deletion, the "Baqqalaa Tolasaa" row has shifted up from
# row 7 to row 6. Remove it as well.
$ws.Rows.Item(6).Delete()

# Leave the selection where Excel naturally ends up after deleting that row.
$ws.Range("A6:XFD6").Select()
